# Adds new Google-Forms response rows (12-15) to the "Form responses 1" sheet,
# matching the cell formatting (style) already used by the existing response rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row used as the formatting template for every new response row (copies its per-cell style).
$templateRow = 11

$newRows = @(
    @{ Row = 12; Cells = @{ "A" = "45430.484464548616"; "B" = "b"; "C" = "ba"; "D" = "bbaabbabb"; "E" = "aa"; "F" = "aa"; "G" = ")("; "H" = "())("; "I" = ")("; "J" = "[]())("; "K" = ")("; "L" = "a"; "M" = "a"; "N" = "aa"; "O" = "aaa"; "P" = "baa" } },
    @{ Row = 13; Cells = @{ "A" = "45433.56760267361"; "B" = "b"; "C" = "b a"; "D" = "a b b"; "F" = "a a"; "G" = ")("; "H" = ")("; "I" = ")("; "J" = ")("; "K" = ")("; "L" = "b"; "M" = "b b"; "P" = "b" } },
    @{ Row = 14; Cells = @{ "A" = "45433.58018291667"; "B" = "aab"; "C" = "ba"; "D" = "bb"; "G" = "][()"; "H" = "]["; "I" = "["; "J" = "][[]"; "K" = "()"; "L" = "bab"; "P" = "bbab" } },
    @{ Row = 15; Cells = @{ "A" = "45433.688659340274"; "C" = "aaba, aaab"; "D" = "bbbbaa, bbbbbbbbaaaa"; "E" = "baabbaab, baabbaabbaab, baabbaabbaabbaab"; "G" = ")[]("; "H" = "(([])), [[]]"; "I" = ")(][)(, [()()]"; "J" = ")()()("; "L" = "abba, ababbaa"; "M" = "baab, abbaab"; "N" = "babaab"; "P" = "bbab" } }
)

foreach ($rowInfo in $newRows) {
    $r = $rowInfo.Row
    foreach ($col in $rowInfo.Cells.Keys) {
        $srcAddr = "$col$templateRow"
        $dstAddr = "$col$r"
        # Copy number/text formatting (style) from the template row so the new cell
        # gets the same style index as the rest of the column (date style for A, text style elsewhere).
        $ws.Range($srcAddr).Copy()
        $ws.Range($dstAddr).PasteSpecial(-4122)
        if ($col -eq "A") {
            $ws.Range($dstAddr).Value = [double]$rowInfo.Cells[$col]
        } else {
            $ws.Range($dstAddr).Value = $rowInfo.Cells[$col]
        }
    }
}